# The original commit bumped the content of the 2x3 data block three times
# in a row, each pass appending "_modified" to whatever was already there:
#   aa/bb/cc/2aa/2bb/2cc
#     -> aa_modified/... (pass 1)
#     -> aa_modified_modified/... (pass 2)
#     -> aa_modified_modified_modified/... (pass 3, final state)
# Re-create that same sequence of edits here so the workbook ends up with
# the same final cell values (and the selection left on H13, like the
# author's last action in Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$addresses = @("A1", "B1", "C1", "A2", "B2", "C2")

foreach ($pass in 1..3) {
    foreach ($addr in $addresses) {
        $current = $ws.Range($addr).Value2
        $ws.Range($addr).Value = "$current" + "_modified"
    }
}

# Leave the selection where the author left it.
$ws.Range("H13").Select()
